# Apply the "simulator full-month coverage, persist logs, fix employees" edit:
#  - Weekly Timesheet: update Client names for rows 2-6
#  - Weekly Timesheet & Jason Schema: populate Rate / Total columns (were 0)
#  - Weekly Timesheet: update subtotal / grand total rows to reflect new totals
#  - Jason Schema: fix Employee ID value

$wb = $excel.ActiveWorkbook

$wsTime = $wb.Worksheets.Item("Weekly Timesheet")
$wsSchema = $wb.Worksheets.Item("Jason Schema")

# --- Weekly Timesheet sheet ---

# Client name corrections (column B, rows 2-6)
$wsTime.Range("B2").Value = "McClure"
$wsTime.Range("B3").Value = "Evans"
$wsTime.Range("B4").Value = "Fritts"
$wsTime.Range("B5").Value = "Hendricks"
$wsTime.Range("B6").Value = "Regan"

# Rate (column E) and Total (column F) for rows 2-6
for ($r = 2; $r -le 6; $r++) {
    $wsTime.Cells.Item($r, 5).Value = 95
    $wsTime.Cells.Item($r, 6).Value = 760
}

# Subtotal / grand total rows
$wsTime.Range("F8").Value = 3800
$wsTime.Range("F11").Value = 3800
$wsTime.Range("F13").Value = 3800

# --- Jason Schema sheet ---

# Client name corrections (column D, rows 2-6) - mirrors Weekly Timesheet
$wsSchema.Range("D2").Value = "McClure"
$wsSchema.Range("D3").Value = "Evans"
$wsSchema.Range("D4").Value = "Fritts"
$wsSchema.Range("D5").Value = "Hendricks"
$wsSchema.Range("D6").Value = "Regan"

# Rate (column F) and Total (column G) for rows 2-6
for ($r = 2; $r -le 6; $r++) {
    $wsSchema.Cells.Item($r, 6).Value = 95
    $wsSchema.Cells.Item($r, 7).Value = 760
}

# Fix Employee ID (every row shares the same Employee ID value)
for ($r = 2; $r -le 6; $r++) {
    $wsSchema.Cells.Item($r, 2).Value = "emp_4nlnrvy7"
}
